$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("P1").Value = "Thurs 2025-05-29"
$ws.Columns("B:O").ColumnWidth = 16.08984375
$ws.Range("B2:N32").HorizontalAlignment = -4131
